# Fruta / hortaliza, semanal
# Insert a new weekly price-report block (3 rows) for "Super Queen" Nectarín
# right before the existing row 276, pushing all subsequent rows down by 3
# (old row 276 -> new row 279, ..., old row 366 -> new row 369).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at 276:278 (Excel shifts rows 276..366 down to 279..369,
# copying formatting - including the date-format style on column D - from row 276).
$ws.Rows("276:278").Insert()

# --- Row 276: Super Queen / Especial ---
$ws.Range("A276").Value2 = 11
$ws.Range("B276").Value = "Vega Monumental Concepción"
$ws.Range("C276").Value = "Bíobío"
$ws.Range("D276").Value2 = 44917
$ws.Range("E276").Value2 = 8
$ws.Range("F276").Value = "Fruta"
$ws.Range("G276").Value2 = 100103
$ws.Range("H276").Value = "Frutos de hueso (carozo)"
$ws.Range("I276").Value2 = 100103006
$ws.Range("J276").Value = "Nectarín"
$ws.Range("K276").Value = "Super Queen"
$ws.Range("L276").Value = "Especial"
$ws.Range("M276").Value2 = 50
$ws.Range("N276").Value2 = 15000
$ws.Range("O276").Value2 = 15000
$ws.Range("P276").Value2 = 15000
$ws.Range("Q276").Value = "$/caja 15 kilos empedrada"
$ws.Range("R276").Value = "Región de O'Higgins"
$ws.Range("S276").Value2 = 1000
$ws.Range("T276").Value2 = 15

# --- Row 277: Super Queen / Primera ---
$ws.Range("A277").Value2 = 11
$ws.Range("B277").Value = "Vega Monumental Concepción"
$ws.Range("C277").Value = "Bíobío"
$ws.Range("D277").Value2 = 44917
$ws.Range("E277").Value2 = 8
$ws.Range("F277").Value = "Fruta"
$ws.Range("G277").Value2 = 100103
$ws.Range("H277").Value = "Frutos de hueso (carozo)"
$ws.Range("I277").Value2 = 100103006
$ws.Range("J277").Value = "Nectarín"
$ws.Range("K277").Value = "Super Queen"
$ws.Range("L277").Value = "Primera"
$ws.Range("M277").Value2 = 100
$ws.Range("N277").Value2 = 13000
$ws.Range("O277").Value2 = 13000
$ws.Range("P277").Value2 = 13000
$ws.Range("Q277").Value = "$/caja 15 kilos empedrada"
$ws.Range("R277").Value = "Región de O'Higgins"
$ws.Range("S277").Value2 = 867
$ws.Range("T277").Value2 = 15

# --- Row 278: Super Queen / Segunda ---
$ws.Range("A278").Value2 = 11
$ws.Range("B278").Value = "Vega Monumental Concepción"
$ws.Range("C278").Value = "Bíobío"
$ws.Range("D278").Value2 = 44917
$ws.Range("E278").Value2 = 8
$ws.Range("F278").Value = "Fruta"
$ws.Range("G278").Value2 = 100103
$ws.Range("H278").Value = "Frutos de hueso (carozo)"
$ws.Range("I278").Value2 = 100103006
$ws.Range("J278").Value = "Nectarín"
$ws.Range("K278").Value = "Super Queen"
$ws.Range("L278").Value = "Segunda"
$ws.Range("M278").Value2 = 50
$ws.Range("N278").Value2 = 11000
$ws.Range("O278").Value2 = 11000
$ws.Range("P278").Value2 = 11000
$ws.Range("Q278").Value = "$/caja 15 kilos empedrada"
$ws.Range("R278").Value = "Región de O'Higgins"
$ws.Range("S278").Value2 = 733
$ws.Range("T278").Value2 = 15
